$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.701.82"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "'1.656.42"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").Value = "'1.003"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").Value = "'303.58"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("E7").Value = "  +0.68%  "
$ws.Range("D8").Value = "'0.3621"
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("D9").Value = "'51.15"
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("D10").Value = "'1.243"
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D11").Value = "'0.08205"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").Value = "'22.66"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").Value = "'6.513"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "'7.440"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").Value = "'0.00001234"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("D17").Value = "'1.655.05"
$ws.Range("E17").Value = "  +1.45%  "
$ws.Range("D18").Value = "'97.49"
$ws.Range("E18").Value = "  +2.32%  "
$ws.Range("D19").Value = "'0.06995"
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("D20").Value = "'6.815"
$ws.Range("E20").Value = "  +3.15%  "
$ws.Range("D21").Value = "'17.70"
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "'12.89"
$ws.Range("E23").Value = "  +2.65%  "
$ws.Range("D24").Value = "'23.698.61"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").Value = "'3.054"
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").Value = "'21.27"
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("D28").Value = "'153.21"
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("D29").Value = "'5.193"
$ws.Range("E29").Value = "  -1.54%  "
$ws.Range("D30").Value = "'134.33"
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("D31").Value = "'1.841.25"
$ws.Range("E31").Value = "  +1.36%  "
$ws.Range("D32").Value = "'6.980"
$ws.Range("E32").Value = "  +5.21%  "
$ws.Range("D33").Value = "'2.197"
$ws.Range("E33").Value = "  +0.96%  "
$ws.Range("D34").Value = "'1.066"
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("D35").Value = "'11.75"
$ws.Range("E35").Value = "  +2.64%  "
$ws.Range("D36").Value = "'0.02822"
$ws.Range("E36").Value = "  +1.96%  "
$ws.Range("D37").Value = "'0.2525"
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("D38").Value = "'6.129"
$ws.Range("E38").Value = "  +1.58%  "
$ws.Range("D39").Value = "'0.08789"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "'0.07065"
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("D41").Value = "'13.04"
$ws.Range("E41").Value = "  +6.91%  "
$ws.Range("D42").Value = "'0.7041"
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("D43").Value = "'1.332"
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("D44").Value = "'15.98"
$ws.Range("E44").Value = "  +1.02%  "
$ws.Range("D45").Value = "'0.6536"
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").Value = "'2.319"
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("D47").Value = "'1.002"
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("D48").Value = "'3.978"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").Value = "'0.07947"
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("D50").Value = "'128.31"
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("D51").Value = "'1.187"
$ws.Range("E51").Value = "  -0.91%  "
